$wb = $excel.ActiveWorkbook

# --- Sheet "建設" (index 15, rId15/sheet15.xml) ---------------------------
# Add the 20 newly scraped article URLs beneath the existing header rows
# (A1 = "Url", A2 = "url_articles"), then select it as the active sheet
# with cell G34 selected -- matching the state Excel saved after the
# "finished get content and download page source" scrape run.
$wsConstruction = $wb.Worksheets.Item(15)

$urls = @(
    "https://ainow.ai/2019/05/23/170454/",
    "https://ainow.ai/2019/05/21/170262/",
    "https://ainow.ai/2019/05/20/170151/",
    "https://ainow.ai/2019/05/12/169561/",
    "https://ainow.ai/2019/05/06/169056/",
    "https://ainow.ai/2019/04/19/168188/",
    "https://ainow.ai/2019/04/14/167755/",
    "https://ainow.ai/2019/04/10/167574/",
    "https://ainow.ai/2019/03/17/165762/",
    "https://ainow.ai/2019/03/06/165166/",
    "https://ainow.ai/2019/02/26/164436/",
    "https://ainow.ai/2019/02/24/164207/",
    "https://ainow.ai/2019/02/19/163875/",
    "https://ainow.ai/2019/01/28/161964/",
    "https://ainow.ai/2019/01/16/160846/",
    "https://ainow.ai/2019/01/15/160768/",
    "https://ainow.ai/2019/01/08/160263/",
    "https://ainow.ai/2019/01/06/159978/",
    "https://ainow.ai/2018/12/27/159554/",
    "https://ainow.ai/2018/12/26/159405/"
)

$row = 3
foreach ($url in $urls) {
    $cell = $wsConstruction.Cells.Item($row, 1)
    $cell.Value = $url
    # The scraped rows were written without inheriting the header column's
    # style (s="1"), so clear formatting back to the workbook default.
    $cell.Style = "Normal"
    $row = $row + 1
}

$wsConstruction.Activate()
$wsConstruction.Range("G34").Select()

# --- Sheet "遠隔会議" (index 8, rId8/sheet8.xml) ---------------------------
# Restore its saved selection (A3:XFD3) without changing which sheet is
# active.
$wsRemoteMeeting = $wb.Worksheets.Item(8)
$wsRemoteMeeting.Range("A3:XFD3").Select()

# Re-activate "建設" so it remains the tab that is selected/visible when the
# workbook is saved (the previous selection on another sheet would otherwise
# steal tabSelected).
$wsConstruction.Activate()
$wsConstruction.Range("G34").Select()
